{"js": "const body = context.document.body;\n\n// --- 1) Merge the \"Word Bias.\" paragraph with the following paragraph ---\n// (they become one paragraph / one flow of text) and extend the sentence\n// describing the Word Bias phenomenon.\nconst wordBiasResults = body.search(\"Word Bias.\", { matchCase: true });\nwordBiasResults.load(\"paragraphs\");\nawait context.sync();\n\nconst wordBiasParagraph = wordBiasResults.items[0].paragraphs.getFirst();\nconst nextParagraph = wordBiasParagraph.getNext();\n\n// Delete the paragraph mark that separates the two paragraphs, merging them\n// into a single paragraph (same effect as pressing Delete at the end of the\n// first paragraph in the Word UI).\nconst paragraphBreak = wordBiasParagraph\n  .getRange(\"End\")\n  .expandTo(nextParagraph.getRange(\"Start\"));\nparagraphBreak.delete();\nawait context.sync();\n\n// Replace the trailing period after \"Word Bias\" with the extended clause.\nconst wordBiasPeriod = body.search(\"Word Bias.\", { matchCase: true });\nawait context.sync();\nwordBiasPeriod.items[0].insertText(\n  \"Word Bias, a phenomenon that affects popular embedding methods such as word2vec and Glove which leads to undesirable word associations ( for instance gender bias). \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 2) \"our corpus of tokens\" -> \"our set of tokens\" ---\nconst corpusResults = body.search(\"our corpus of tokens\", { matchCase: true });\nawait context.sync();\ncorpusResults.items[0].insertText(\"our set of tokens\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 3) Append a new closing sentence after \"pronoun.\" ---\nconst pronounResults = body.search(\"pronoun.\", { matchCase: true });\nawait context.sync();\npronounResults.items[0].insertText(\n  \" As we hoped, removing a small part of the training corpus led to an improvement in the overall bias.\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) Merge the \"Word Bias.\" paragraph with the following paragraph and\n#        extend the sentence describing the Word Bias phenomenon. ---\n\n# Locate the paragraph that ends in \"... Word Bias.\" by scanning the\n# document's Paragraphs collection (more reliable than deriving a paragraph\n# from a small Find match range).\n$paras = $d.Paragraphs\n$wordBiasPara = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"*Word Bias.*\") {\n        $wordBiasPara = $p\n        break\n    }\n}\n\n# Delete the paragraph mark at the end of that paragraph so it merges with\n# the following paragraph (same as pressing Delete at the end of the\n# paragraph, right before \"This problem cannot be fixed...\", in the Word UI).\n$pRange = $wordBiasPara.Range\n$mark = $d.Range($pRange.End - 1, $pRange.End)\n$mark.Delete()\n\n# Replace the trailing period after \"Word Bias\" with the extended clause.\n$rng = $d.Content\n$rng.Find.Execute(\"Word Bias.\")\n$rng.Text = \"Word Bias, a phenomenon that affects popular embedding methods such as word2vec and Glove which leads to undesirable word associations ( for instance gender bias). \"\n\n# --- 2) \"our corpus of tokens\" -> \"our set of tokens\" ---\n$rng = $d.Content\n$rng.Find.Execute(\"our corpus of tokens\")\n$rng.Text = \"our set of tokens\"\n\n# --- 3) Append a new closing sentence after \"pronoun.\" ---\n$rng = $d.Content\n$rng.Find.Execute(\"pronoun.\")\n$rng.Text = \"pronoun. As we hoped, removing a small part of the training corpus led to an improvement in the overall bias.\"\n"}
